# Powerpoint writer: avoid extra blank lines before author.
# (In the case where there is no subtitle.)
#
# The "Testing Layouts" title slide has an empty Subtitle placeholder
# (just two manual line breaks) left over on it; remove it entirely so
# no blank subtitle shape remains on the slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("Subtitle 2").Cut()
